$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.910.03'
$ws.Range("E2").Value = '  +6.42%  '

$ws.Range("D3").Value = '3.114.74'
$ws.Range("E3").Value = '  +5.18%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +13.31%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").Value = '3.110.04'
$ws.Range("E8").Value = '  +5.30%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.501'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +19.77%  '

$ws.Range("E11").Value = '  +8.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.04%  '

$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.52%  '

$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000228'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.95%  '

$ws.Range("D15").Value = '3.595.47'
$ws.Range("E15").Value = '  +4.37%  '

$ws.Range("D16").Value = '64.945.37'
$ws.Range("E16").Value = '  +6.39%  '

$ws.Range("D17").Value = '3.104.74'
$ws.Range("E17").Value = '  +4.79%  '

$ws.Range("E18").Value = '  -0.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '488.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.676'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +12.96%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.21'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.87%  '

$ws.Range("E26").Value = '  +0.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.83'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.04%  '

$ws.Range("E29").Value = '  +11.11%  '

$ws.Range("E30").Value = '  -0.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.25'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.65%  '

$ws.Range("E32").Value = '  +4.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.46'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.95%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.78'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.65'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.15'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.88%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '469.54'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.60%  '

$ws.Range("E38").Value = '  +9.99%  '

$ws.Range("E39").Value = '  +7.00%  '

$ws.Range("D40").Value = '3.037.85'
$ws.Range("E40").Value = '  -1.35%  '

$ws.Range("E41").Value = '  +3.54%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +19.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '28.38'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +14.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.262'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.16%  '

$ws.Range("E46").Value = '  +0.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.90%  '

$ws.Range("E48").Value = '  +6.17%  '

$ws.Range("D49").Value = '0.0₃0522'
$ws.Range("E49").Value = '  +11.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '118.24'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.08'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.95%  '
